$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '34.136.38'
$ws.Range("E2").Value = '  -1.44%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.785.04'
$ws.Range("E3").Value = '  -1.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '221.71'
$ws.Range("E5").Value = '  -1.84%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.550'
$ws.Range("E6").Value = '  -1.11%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  +0.12%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '31.65'
$ws.Range("E8").Value = '  -3.96%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.289'
$ws.Range("E9").Value = '  +1.36%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0710'
$ws.Range("E10").Value = '  +6.45%  '

$ws.Range("E11").Value = '  -1.03%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.041.51'
$ws.Range("E12").Value = '  -0.93%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.790.58'
$ws.Range("E13").Value = '  -0.86%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '10.71'
$ws.Range("E14").Value = '  -3.66%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.626'
$ws.Range("E15").Value = '  -2.62%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '34.091.51'
$ws.Range("E16").Value = '  -1.41%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.23'
$ws.Range("E17").Value = '  -1.31%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '68.08'
$ws.Range("E18").Value = '  -2.37%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '245.06'
$ws.Range("E19").Value = '  -4.47%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0779'
$ws.Range("E20").Value = '  +2.73%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.999'
$ws.Range("E21").Value = '  +0.18%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.67'
$ws.Range("E22").Value = '  +1.73%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.11'
$ws.Range("E23").Value = '  -3.18%  '

$ws.Range("E24").Value = '  -0.60%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '157.71'
$ws.Range("E25").Value = '  -0.14%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '16.37'
$ws.Range("E26").Value = '  -0.85%  '

$ws.Range("E27").Value = '  -1.43%  '

$ws.Range("E28").Value = '  -1.76%  '

$ws.Range("E29").Value = '  +0.08%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0520'
$ws.Range("E30").Value = '  +0.43%  '

$ws.Range("E31").Value = '  -2.61%  '

$ws.Range("E32").Value = '  +0.75%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.53'
$ws.Range("E33").Value = '  -2.08%  '

$ws.Range("E34").Value = '  -4.19%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.394.69'
$ws.Range("E35").Value = '  -4.72%  '

$ws.Range("E36").Value = '  -0.66%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.629'
$ws.Range("E37").Value = '  -0.88%  '

$ws.Range("E38").Value = '  -1.87%  '

$ws.Range("B39").Value = 'ARBITRUM'
$ws.Range("C39").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.935'
$ws.Range("E39").Value = '  +3.76%  '

$ws.Range("B40").Value = 'Aave'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '79.63'
$ws.Range("E40").Value = '  -4.54%  '

$ws.Range("B41").Value = 'HuobiToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.34'
$ws.Range("E41").Value = '  +1.20%  '

$ws.Range("B42").Value = 'MXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.71'
$ws.Range("E42").Value = '  -5.28%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.12'
$ws.Range("E43").Value = '  +1.43%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0493'
$ws.Range("E44").Value = '  -3.01%  '

$ws.Range("B45").Value = 'FraxShare'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.84'
$ws.Range("E45").Value = '  -1.27%  '

$ws.Range("B46").Value = 'WEMIXToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.04'
$ws.Range("E46").Value = '  +0.05%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '106.09'
$ws.Range("E47").Value = '  +5.31%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.937.61'
$ws.Range("E48").Value = '  -1.13%  '

$ws.Range("E49").Value = '  -0.30%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '11.83'
$ws.Range("E50").Value = '  -0.54%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0₆0119'
$ws.Range("E51").Value = '  +2.25%  '
